{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change (per the supplied diff):\n//   1. The first paragraph (\"Status  hello  how are you\") gets Word's\n//      automatic grammar-check markers (<w:proofErr w:type=\"gramStart\"/>\n//      ... <w:proofErr w:type=\"gramEnd\"/>) wrapped around \"Status  hello\",\n//      which also forces that text to be split across an extra run\n//      boundary (no visible text change).\n//   2. A brand-new second paragraph \"Hi iam feature branch\" is added,\n//      with Word's spell-check markers (<w:proofErr w:type=\"spellStart\"/>\n//      ... <w:proofErr w:type=\"spellEnd\"/>) wrapped around the\n//      not-in-dictionary word \"iam\".\n//   3. The \"_GoBack\" bookmark (Word's \"last edit position\" marker) ends\n//      up at the end of the new second paragraph instead of the first.\n//\n// `w:proofErr` is not part of the Word JS object model, so the only way\n// to reproduce it faithfully is to splice literal OOXML for the affected\n// paragraphs via insertOoxml(..., \"Replace\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that holds the original \"Status ... hello ...\" text\n// (falls back to the first paragraph if nothing matches).\nlet target = paragraphs.items.find(p => p.text.indexOf(\"Status\") !== -1 && p.text.indexOf(\"hello\") !== -1);\nif (!target) {\n  target = paragraphs.items[0];\n}\n\n// Step 1: add the new second paragraph right after it (plain text for now;\n// its final OOXML - including the spell-check proof markers - is applied\n// below).\nconst newPara = target.insertParagraph(\"Hi iam feature branch\", \"After\");\nawait context.sync();\n\nconst WORD_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\nfunction wordPackage(paragraphXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + WORD_NS + '><w:body>' + paragraphXml + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>';\n}\n\n// Step 2: rewrite paragraph 1 with the gramStart/gramEnd proof-error split.\nconst para1Xml =\n  '<w:p>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">Status </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> hello</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">  how are you</w:t></w:r>' +\n  '</w:p>';\ntarget.getRange(\"Whole\").insertOoxml(wordPackage(para1Xml), \"Replace\");\nawait context.sync();\n\n// Step 3: rewrite paragraph 2 with the spellStart/spellEnd proof-error\n// around \"iam\", and carry the \"_GoBack\" bookmark onto its end.\nconst para2Xml =\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Hi </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>iam</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> feature branch</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\nnewPara.getRange(\"Whole\").insertOoxml(wordPackage(para2Xml), \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) edit script.\n#\n# Target change (per the supplied diff):\n#   1. The first paragraph (\"Status  hello  how are you\") gets Word's\n#      automatic grammar-check markers (<w:proofErr w:type=\"gramStart\"/>\n#      ... <w:proofErr w:type=\"gramEnd\"/>) wrapped around \"Status  hello\",\n#      which also forces that text to be split across an extra run\n#      boundary (no visible text change).\n#   2. A brand-new second paragraph \"Hi iam feature branch\" is added,\n#      with Word's spell-check markers (<w:proofErr w:type=\"spellStart\"/>\n#      ... <w:proofErr w:type=\"spellEnd\"/>) wrapped around the\n#      not-in-dictionary word \"iam\".\n#   3. The \"_GoBack\" bookmark (Word's \"last edit position\" marker) ends\n#      up at the end of the new second paragraph instead of the first.\n#\n# `w:proofErr` has no dedicated COM property/method, so the only faithful\n# way to reproduce it is to splice literal OOXML for the affected\n# paragraphs via Range.InsertXML (the COM counterpart of Office.js's\n# insertOoxml), using the same pkg:package wrapper Range.WordOpenXML\n# returns.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the original \"Status ... hello ...\" text\n# (falls back to the first paragraph if nothing matches).\n$targetIndex = 1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t -like \"*Status*\" -and $t -like \"*hello*\") {\n    $targetIndex = $i\n    break\n  }\n}\n\n# Step 1: detach the \"_GoBack\" bookmark from its current spot; it will be\n# re-added (explicitly, in OOXML) at the end of the new second paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 2: add the new second paragraph right after the target paragraph\n# (placeholder text for now; its final OOXML - including the spell-check\n# proof markers and the bookmark - is applied below).\n$p1 = $d.Paragraphs.Item($targetIndex)\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($targetIndex + 1)\n$p2.Range.Text = \"Hi iam feature branch\"\n\nfunction Wrap-WordOoxml($bodyInner) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData>' + `\n    '</pkg:part></pkg:package>'\n}\n\n# Step 3: rewrite paragraph 1 with the gramStart/gramEnd proof-error split.\n$p1 = $d.Paragraphs.Item($targetIndex)\n$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)\n$para1Xml = '<w:p>' + `\n    '<w:proofErr w:type=\"gramStart\"/>' + `\n    '<w:r><w:t xml:space=\"preserve\">Status </w:t></w:r>' + `\n    '<w:r><w:t xml:space=\"preserve\"> hello</w:t></w:r>' + `\n    '<w:proofErr w:type=\"gramEnd\"/>' + `\n    '<w:r><w:t xml:space=\"preserve\">  how are you</w:t></w:r>' + `\n  '</w:p>'\n$r1.InsertXML((Wrap-WordOoxml $para1Xml))\n\n# Step 4: rewrite paragraph 2 with the spellStart/spellEnd proof-error\n# around \"iam\", and carry the \"_GoBack\" bookmark onto its end.\n$p2 = $d.Paragraphs.Item($targetIndex + 1)\n$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)\n$para2Xml = '<w:p>' + `\n    '<w:r><w:t xml:space=\"preserve\">Hi </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:t>iam</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:t xml:space=\"preserve\"> feature branch</w:t></w:r>' + `\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' + `\n    '<w:bookmarkEnd w:id=\"0\"/>' + `\n  '</w:p>'\n$r2.InsertXML((Wrap-WordOoxml $para2Xml))\n"}
